$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 6.55
$ws.Range("Q22").Value = 2.2
$ws.Range("R22").Value = 1.65
$ws.Range("G24").Value = 2.1
$ws.Range("H24").Value = 3.2
$ws.Range("I24").Value = 3.75
$ws.Range("L24").Value = 4
$ws.Range("N24").Value = 8.5
$ws.Range("Q24").Value = 2.08
$ws.Range("R24").Value = 1.73
$ws.Range("W24").Value = 7.5
$ws.Range("Z24").Value = 19
$ws.Range("AC24").Value = 8.5
$ws.Range("AK24").Value = 29
$ws.Range("BA24").Value = 81
$ws.Range("J27").Value = 2.82
$ws.Range("K27").Value = 2.1
$ws.Range("L27").Value = 3.55
$ws.Range("N27").Value = 7
$ws.Range("O27").Value = 1.32
$ws.Range("P27").Value = 3.1
$ws.Range("S27").Value = 1.4
$ws.Range("T27").Value = 2.75
$ws.Range("V27").Value = 1.98
$ws.Range("W27").Value = 7.5
$ws.Range("X27").Value = 10.75
$ws.Range("Y27").Value = 9
$ws.Range("AA27").Value = 18.5
$ws.Range("AB27").Value = 29
$ws.Range("AC27").Value = 7
$ws.Range("AE27").Value = 13.5
$ws.Range("AG27").Value = 9.5
$ws.Range("AH27").Value = 16
$ws.Range("AI27").Value = 10.75
$ws.Range("AK27").Value = 26
$ws.Range("AL27").Value = 32
$ws.Range("AO27").Value = 11.75
$ws.Range("AP27").Value = 19.5
$ws.Range("AR27").Value = 80
$ws.Range("AS27").Value = 250
$ws.Range("AT27").Value = 2.75
$ws.Range("AU27").Value = 6.9
$ws.Range("AV27").Value = 60
$ws.Range("AW27").Value = 5
$ws.Range("AX27").Value = 16.5
$ws.Range("AZ27").Value = 75
$ws.Range("G35").Value = 1.8
$ws.Range("H35").Value = 3.35
$ws.Range("I35").Value = 4.2
$ws.Range("J35").Value = 2.35
$ws.Range("K35").Value = 2.12
$ws.Range("L35").Value = 4.55
$ws.Range("O35").Value = 1.31
$ws.Range("P35").Value = 2.87
$ws.Range("Q35").Value = 1.98
$ws.Range("U35").Value = 1.8
$ws.Range("V35").Value = 1.8
$ws.Range("W35").Value = 6.6
$ws.Range("AE35").Value = 15.5
$ws.Range("AF35").Value = 75
$ws.Range("AG35").Value = 11
$ws.Range("AH35").Value = 23
$ws.Range("AJ35").Value = 70
$ws.Range("AK35").Value = 45
$ws.Range("AL35").Value = 50
$ws.Range("AM35").Value = 700
$ws.Range("AO35").Value = 8.75
$ws.Range("AP35").Value = 17
$ws.Range("AQ35").Value = 30
$ws.Range("AU35").Value = 7.1
$ws.Range("AV35").Value = 65
$ws.Range("AX35").Value = 24
$ws.Range("AY35").Value = 29
$ws.Range("AZ35").Value = 150
$ws.Range("BA35").Value = 175
$ws.Range("BB35").Value = 350
